# Implement quantitative metrics highlighting (bold + color #2C3E50)
# for impact metrics (percentages, dollar amounts, large numbers) across
# the achievements / work-experience bullet points of the resume.

$d = $word.ActiveDocument
$HIGHLIGHT_COLOR = 5258796   # RGB(44, 62, 80) -> hex 2C3E50

# Finds $searchText inside $paragraph (a Word Paragraph object) and, if
# found, applies bold + the highlight color to just that sub-range. Word
# automatically splits the run(s) so the rest of the paragraph's text
# keeps its original (unformatted) run.
function Set-MetricHighlight {
    param(
        $Paragraph,
        [string]$SearchText
    )

    $rng = $Paragraph.Range
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $HIGHLIGHT_COLOR
    }
    return $found
}

# Locates the single paragraph whose text contains $containsText (a unique
# identifying substring of the whole, not-yet-modified paragraph) -- and
# which does NOT contain $excludeText (used to disambiguate paragraphs
# that otherwise share a common prefix) -- and applies Set-MetricHighlight
# for every metric string supplied afterwards, left to right, matching the
# order they occur in the text.
function Format-ParagraphMetrics {
    param(
        [string]$ContainsText,
        [string[]]$Metrics,
        [string]$ExcludeText = $null
    )

    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t -like "*$ContainsText*") {
            if ($ExcludeText -and ($t -like "*$ExcludeText*")) {
                continue
            }
            foreach ($metric in $Metrics) {
                Set-MetricHighlight $p $metric | Out-Null
            }
            return
        }
    }
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
Format-ParagraphMetrics "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial" @("23%", "64%")

# 2) "Achieved 87% prediction accuracy ... ±4.2% to ±2.1%" (long version, under Siege Analytics)
Format-ParagraphMetrics "reducing polling error margins" @("87%", "71%", [char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%")

# 3) "Wrote RFP and analyzed bids from 1,200 vendors ..."
Format-ParagraphMetrics "Wrote RFP and analyzed bids from" @("1,200")

# 4) "Created comprehensive meta-analysis framework ... $400M ... $1B+"
Format-ParagraphMetrics "Created comprehensive meta-analysis framework" @("`$400M", "`$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Format-ParagraphMetrics "Algorithm reduced mapping costs by" @("73.5%", "`$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short version, Key Achievements)
Format-ParagraphMetrics "Achieved" @("87%", "71%") "reducing polling error margins"
